# Auto-generated edit script for NORTH_CAROLINA_2017.xlsx-style workbook edit
# 1) Rename header row columns to snake_case analytic names
# 2) Title-case Spanish linking particles (de/del/el/la/las/los/y) in place names
# 3) Bump 7 D-column percentage values by one ULP (0.0009340390517279723 -> ...724)
# 4) Remove the trailing footnote/metadata rows (1596-1601), shrinking the used range to A1:D1595

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row rename ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Title-case the "de/del/el/la/las/los/y" particles in place names ---
$ws.Range("B7").Value = "Pabellón De Arteaga"
$ws.Range("B8").Value = "Rincón De Romos"
$ws.Range("B31").Value = "Amatenango De La Frontera"
$ws.Range("B32").Value = "Amatenango Del Valle"
$ws.Range("B36").Value = "Benemérito De Las Américas"
$ws.Range("B44").Value = "Chiapa De Corzo"
$ws.Range("B49").Value = "Comitán De Domínguez"
$ws.Range("B71").Value = "Mazapa De Madero"
$ws.Range("B76").Value = "Ocozocoautla De Espinosa"
$ws.Range("B85").Value = "Salto De Agua"
$ws.Range("B86").Value = "San Cristóbal De Las Casas"
$ws.Range("B124").Value = "Hidalgo Del Parral"
$ws.Range("B134").Value = "San Francisco Del Oro"
$ws.Range("B153").Value = "San Juan De Sabinas"
$ws.Range("A165").Value = "Ciudad De México"
$ws.Range("B169").Value = "Cuajimalpa De Morelos"
$ws.Range("B183").Value = "Coneto De Comonfort"
$ws.Range("B196").Value = "Nombre De Dios"
$ws.Range("B200").Value = "Pánuco De Coronado"
$ws.Range("B206").Value = "San Juan Del Río"
$ws.Range("A216").Value = "Estado De México"
$ws.Range("B216").Value = "Acambay De Ruíz Castañeda"
$ws.Range("B219").Value = "Almoloya De Alquisiras"
$ws.Range("B220").Value = "Almoloya De Juárez"
$ws.Range("B227").Value = "Atizapán De Zaragoza"
$ws.Range("B237").Value = "Coacalco De Berriozábal"
$ws.Range("B243").Value = "Ecatepec De Morelos"
$ws.Range("B249").Value = "Ixtapan De La Sal"
$ws.Range("B262").Value = "Naucalpan De Juárez"
$ws.Range("B272").Value = "San Antonio La Isla"
$ws.Range("B273").Value = "San Felipe Del Progreso"
$ws.Range("B275").Value = "San Simón De Guerrero"
$ws.Range("B285").Value = "Tenango Del Valle"
$ws.Range("B296").Value = "Tlalnepantla De Baz"
$ws.Range("B301").Value = "Valle De Bravo"
$ws.Range("B302").Value = "Valle De Chalco Solidaridad"
$ws.Range("B303").Value = "Villa De Allende"
$ws.Range("B304").Value = "Villa Del Carbón"
$ws.Range("B317").Value = "San Miguel De Allende"
$ws.Range("B318").Value = "Apaseo El Alto"
$ws.Range("B319").Value = "Apaseo El Grande"
$ws.Range("B327").Value = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Range("B331").Value = "Jaral Del Progreso"
$ws.Range("B339").Value = "Purísima Del Rincón"
$ws.Range("B343").Value = "San Diego De La Unión"
$ws.Range("B345").Value = "San Francisco Del Rincón"
$ws.Range("B347").Value = "San Luis De La Paz"
$ws.Range("B348").Value = "Santa Cruz De Juventino Rosas"
$ws.Range("B349").Value = "Silao De La Victoria"
$ws.Range("B354").Value = "Valle De Santiago"
$ws.Range("B360").Value = "Acapulco De Juárez"
$ws.Range("B363").Value = "Ajuchitlán Del Progreso"
$ws.Range("B364").Value = "Alcozauca De Guerrero"
$ws.Range("B367").Value = "Atenango Del Río"
$ws.Range("B368").Value = "Atlamajalcingo Del Monte"
$ws.Range("B370").Value = "Atoyac De Álvarez"
$ws.Range("B371").Value = "Ayutla De Los Libres"
$ws.Range("B374").Value = "Buenavista De Cuéllar"
$ws.Range("B375").Value = "Chilapa De Álvarez"
$ws.Range("B376").Value = "Chilpancingo De Los Bravo"
$ws.Range("B380").Value = "Coyuca De Benítez"
$ws.Range("B381").Value = "Coyuca De Catalán"
$ws.Range("B385").Value = "Cuetzala Del Progreso"
$ws.Range("B386").Value = "Cutzamala De Pinzón"
$ws.Range("B392").Value = "Huitzuco De Los Figueroa"
$ws.Range("B393").Value = "Iguala De La Independencia"
$ws.Range("B395").Value = "Ixcateopan De Cuauhtémoc"
$ws.Range("B396").Value = "Zihuatanejo De Azueta"
$ws.Range("B398").Value = "La Unión De Isidoro Montes De Oca"
$ws.Range("B401").Value = "Mártir De Cuilapan"
$ws.Range("B414").Value = "Taxco De Alarcón"
$ws.Range("B416").Value = "Técpan De Galeana"
$ws.Range("B418").Value = "Tepecoacuilco De Trujano"
$ws.Range("B420").Value = "Tixtla De Guerrero"
$ws.Range("B423").Value = "Tlalixtaquilla De Maldonado"
$ws.Range("B424").Value = "Tlapa De Comonfort"
$ws.Range("B436").Value = "Agua Blanca De Iturbide"
$ws.Range("B442").Value = "Atotonilco De Tula"
$ws.Range("B443").Value = "Atotonilco El Grande"
$ws.Range("B447").Value = "Cuautepec De Hinojosa"
$ws.Range("B450").Value = "Huasca De Ocampo"
$ws.Range("B453").Value = "Huejutla De Reyes"
$ws.Range("B456").Value = "Jacala De Ledezma"
$ws.Range("B462").Value = "Mineral Del Chico"
$ws.Range("B463").Value = "Mineral Del Monte"
$ws.Range("B464").Value = "Mixquiahuala De Juárez"
$ws.Range("B465").Value = "Molango De Escamilla"
$ws.Range("B467").Value = "Nopala De Villagrán"
$ws.Range("B468").Value = "Pachuca De Soto"
$ws.Range("B471").Value = "Progreso De Obregón"
$ws.Range("B477").Value = "Santiago De Anaya"
$ws.Range("B478").Value = "Santiago Tulantepec De Lugo Guerrero"
$ws.Range("B482").Value = "Tenango De Doria"
$ws.Range("B484").Value = "Tepehuacán De Guerrero"
$ws.Range("B485").Value = "Tepeji Del Río De Ocampo"
$ws.Range("B488").Value = "Tezontepec De Aldama"
$ws.Range("B496").Value = "Tula De Allende"
$ws.Range("B497").Value = "Tulancingo De Bravo"
$ws.Range("B500").Value = "Zacualtipán De Ángeles"
$ws.Range("B501").Value = "Zapotlán De Juárez"
$ws.Range("B505").Value = "Ahualulco De Mercado"
$ws.Range("B509").Value = "Autlán De Navarro"
$ws.Range("B516").Value = "Concepción De Buenos Aires"
$ws.Range("B523").Value = "Encarnación De Díaz"
$ws.Range("B529").Value = "Huejuquilla El Alto"
$ws.Range("B530").Value = "Ixtlahuacán Del Río"
$ws.Range("B533").Value = "Jilotlán De Los Dolores"
$ws.Range("B537").Value = "Lagos De Moreno"
$ws.Range("B540").Value = "Ojuelos De Jalisco"
$ws.Range("B546").Value = "San Juan De Los Lagos"
$ws.Range("B548").Value = "San Miguel El Alto"
$ws.Range("B549").Value = "San Sebastián Del Oeste"
$ws.Range("B550").Value = "Santa María De Los Ángeles"
$ws.Range("B553").Value = "Tamazula De Gordiano"
$ws.Range("B558").Value = "Tepatitlán De Morelos"
$ws.Range("B561").Value = "Tizapán El Alto"
$ws.Range("B562").Value = "Tlajomulco De Zúñiga"
$ws.Range("B567").Value = "Unión De San Antonio"
$ws.Range("B568").Value = "Unión De Tula"
$ws.Range("B569").Value = "Valle De Guadalupe"
$ws.Range("B574").Value = "Yahualica De González Gallo"
$ws.Range("B575").Value = "Zacoalco De Torres"
$ws.Range("B577").Value = "Zapotitlán De Vadillo"
$ws.Range("B578").Value = "Zapotlán El Grande"
$ws.Range("B601").Value = "Coalcomán De Vázquez Pallares"
$ws.Range("B663").Value = "Tiquicheo De Nicolás Romero"
$ws.Range("B687").Value = "Coatlán Del Río"
$ws.Range("B694").Value = "Jonacatepec De Leandro Valle"
$ws.Range("B698").Value = "Puente De Ixtla"
$ws.Range("B704").Value = "Tetela Del Volcán"
$ws.Range("B705").Value = "Tlaltizapán De Zapata"
$ws.Range("B712").Value = "Zacualpan De Amilpas"
$ws.Range("B716").Value = "Amatlán De Cañas"
$ws.Range("B717").Value = "Bahía De Banderas"
$ws.Range("B719").Value = "Ixtlán Del Río"
$ws.Range("B726").Value = "Santa María Del Oro"
$ws.Range("B747").Value = "Lampazos De Naranjo"
$ws.Range("B749").Value = "Mier Y Noriega"
$ws.Range("B754").Value = "San Nicolás De Los Garza"
$ws.Range("B759").Value = "Acatlán De Pérez Figueroa"
$ws.Range("B766").Value = "Capulálpam De Méndez"
$ws.Range("B768").Value = "Chalcatongo De Hidalgo"
$ws.Range("B769").Value = "Ciénega De Zimatlán"
$ws.Range("B773").Value = "Constancia Del Rosario"
$ws.Range("B775").Value = "Fresnillo De Trujano"
$ws.Range("B776").Value = "Guadalupe De Ramírez"
$ws.Range("B778").Value = "Guelatao De Juárez"
$ws.Range("B779").Value = "Guevea De Humboldt"
$ws.Range("B780").Value = "Heroica Ciudad De Ejutla De Crespo"
$ws.Range("B781").Value = "Heroica Ciudad De Huajuapan De León"
$ws.Range("B782").Value = "Heroica Ciudad De Tlaxiaco"
$ws.Range("B783").Value = "Ixtlán De Juárez"
$ws.Range("B784").Value = "Heroica Ciudad De Juchitán De Zaragoza"
$ws.Range("B791").Value = "Mártires De Tacubaya"
$ws.Range("B794").Value = "Miahuatlán De Porfirio Díaz"
$ws.Range("B797").Value = "Nejapa De Madero"
$ws.Range("B798").Value = "Oaxaca De Juárez"
$ws.Range("B799").Value = "Ocotlán De Morelos"
$ws.Range("B800").Value = "Pinotepa De Don Luis"
$ws.Range("B802").Value = "Putla Villa De Guerrero"
$ws.Range("B812").Value = "San Antonio De La Cal"
$ws.Range("B822").Value = "San Dionisio Del Mar"
$ws.Range("B825").Value = "San Felipe Jalapa De Díaz"
$ws.Range("B827").Value = "San Francisco Del Mar"
$ws.Range("B838").Value = "San José Del Peñasco"
$ws.Range("B844").Value = "San Juan Bautista Lo De Soto"
$ws.Range("B886").Value = "San Miguel Del Puerto"
$ws.Range("B887").Value = "San Miguel El Grande"
$ws.Range("B900").Value = "San Pablo Villa De Mitla"
$ws.Range("B903").Value = "San Pedro El Alto"
$ws.Range("B928").Value = "Santa Ana Del Valle"
$ws.Range("B939").Value = "Santa Cruz De Bravo"
$ws.Range("B943").Value = "Santa Cruz Tacache De Mina"
$ws.Range("B949").Value = "Santa Inés Del Monte"
$ws.Range("B959").Value = "Santa María Del Rosario"
$ws.Range("B1012").Value = "Santo Domingo De Morelos"
$ws.Range("B1026").Value = "Tamazulápam Del Espíritu Santo"
$ws.Range("B1027").Value = "Tanetze De Zaragoza"
$ws.Range("B1028").Value = "Tataltepec De Valdés"
$ws.Range("B1029").Value = "Teococuilco De Marcos Pérez"
$ws.Range("B1030").Value = "Teotitlán De Flores Magón"
$ws.Range("B1031").Value = "Heroica Villa Tezoatlán De Segura Y Luna, Cuna De La Independencia De Oaxaca"
$ws.Range("B1032").Value = "Tlacolula De Matamoros"
$ws.Range("B1033").Value = "Tlalixtac De Cabrera"
$ws.Range("B1034").Value = "Totontepec Villa De Morelos"
$ws.Range("B1036").Value = "Villa De Chilapa De Díaz"
$ws.Range("B1037").Value = "Villa De Etla"
$ws.Range("B1038").Value = "Villa De Tamazulápam Del Progreso"
$ws.Range("B1039").Value = "Villa De Tututepec"
$ws.Range("B1040").Value = "Villa De Zaachila"
$ws.Range("B1043").Value = "Villa Sola De Vega"
$ws.Range("B1044").Value = "Villa Talea De Castro"
$ws.Range("B1046").Value = "Zapotitlán Del Río"
$ws.Range("B1048").Value = "Zimatlán De Álvarez"
$ws.Range("B1063").Value = "Ayotoxco De Guerrero"
$ws.Range("B1065").Value = "Chalchicomula De Sesma"
$ws.Range("B1081").Value = "Cuapiaxtla De Madero"
$ws.Range("B1084").Value = "Cuayuca De Andrade"
$ws.Range("B1085").Value = "Cuetzalan Del Progreso"
$ws.Range("B1094").Value = "Huehuetlán El Chico"
$ws.Range("B1095").Value = "Huehuetlán El Grande"
$ws.Range("B1099").Value = "Huitzilan De Serdán"
$ws.Range("B1101").Value = "Ixcamilpa De Guerrero"
$ws.Range("B1103").Value = "Izúcar De Matamoros"
$ws.Range("B1112").Value = "Los Reyes De Juárez"
$ws.Range("B1118").Value = "Palmar De Bravo"
$ws.Range("B1136").Value = "San Salvador El Seco"
$ws.Range("B1137").Value = "San Salvador El Verde"
$ws.Range("B1142").Value = "Tecali De Herrera"
$ws.Range("B1149").Value = "Tepanco De López"
$ws.Range("B1150").Value = "Tepatlaxco De Hidalgo"
$ws.Range("B1155").Value = "Tepexi De Rodríguez"
$ws.Range("B1157").Value = "Tepeyahualco De Cuauhtémoc"
$ws.Range("B1158").Value = "Tetela De Ocampo"
$ws.Range("B1159").Value = "Teteles De Avila Castillo"
$ws.Range("B1164").Value = "Tlacotepec De Benito Juárez"
$ws.Range("B1174").Value = "Totoltepec De Guerrero"
$ws.Range("B1176").Value = "Tuzamapan De Galeana"
$ws.Range("B1197").Value = "Amealco De Bonfil"
$ws.Range("B1199").Value = "Cadereyta De Montes"
$ws.Range("B1205").Value = "Jalpan De Serra"
$ws.Range("B1206").Value = "Landa De Matamoros"
$ws.Range("B1209").Value = "Pinal De Amoles"
$ws.Range("B1212").Value = "San Juan Del Río"
$ws.Range("B1223").Value = "Armadillo De Los Infante"
$ws.Range("B1229").Value = "Ciudad Del Maíz"
$ws.Range("B1238").Value = "Mexquitic De Carmona"
$ws.Range("B1243").Value = "San Ciro De Acosta"
$ws.Range("B1248").Value = "Santa María Del Río"
$ws.Range("B1250").Value = "Soledad De Graciano Sánchez"
$ws.Range("B1260").Value = "Villa De Arista"
$ws.Range("B1261").Value = "Villa De Arriaga"
$ws.Range("B1262").Value = "Villa De Guadalupe"
$ws.Range("B1263").Value = "Villa De Ramos"
$ws.Range("B1264").Value = "Villa De Reyes"
$ws.Range("B1307").Value = "Jalpa De Méndez"
$ws.Range("B1340").Value = "Soto La Marina"
$ws.Range("B1357").Value = "Nanacamilpa De Mariano Arista"
$ws.Range("B1360").Value = "San Pablo Del Monte"
$ws.Range("B1363").Value = "Tepetitla De Lardizábal"
$ws.Range("B1366").Value = "Tetla De La Solidaridad"
$ws.Range("B1383").Value = "Alto Lucero De Gutiérrez Barrios"
$ws.Range("B1387").Value = "Amatlán De Los Reyes"
$ws.Range("B1394").Value = "Boca Del Río"
$ws.Range("B1399").Value = "Castillo De Teayo"
$ws.Range("B1401").Value = "Cazones De Herrera"
$ws.Range("B1420").Value = "Cosamaloapan De Carpio"
$ws.Range("B1421").Value = "Cosautlán De Carvajal"
$ws.Range("B1437").Value = "Hueyapan De Ocampo"
$ws.Range("B1438").Value = "Ignacio De La Llave"
$ws.Range("B1442").Value = "Ixhuatlán De Madero"
$ws.Range("B1443").Value = "Ixhuatlán Del Café"
$ws.Range("B1444").Value = "Ixhuatlán Del Sureste"
$ws.Range("B1453").Value = "Juchique De Ferrer"
$ws.Range("B1456").Value = "Las Vigas De Ramírez"
$ws.Range("B1457").Value = "Lerdo De Tejada"
$ws.Range("B1461").Value = "Martínez De La Torre"
$ws.Range("B1463").Value = "Medellín De Bravo"
$ws.Range("B1467").Value = "Nanchital De Lázaro Cárdenas Del Río"
$ws.Range("B1476").Value = "Ozuluama De Mascareñas"
$ws.Range("B1480").Value = "Paso De Ovejas"
$ws.Range("B1481").Value = "Paso Del Macho"
$ws.Range("B1484").Value = "Poza Rica De Hidalgo"
$ws.Range("B1492").Value = "Sayula De Alemán"
$ws.Range("B1495").Value = "Soledad De Doblado"
$ws.Range("B1500").Value = "Tatahuicapan De Juárez"
$ws.Range("B1530").Value = "Vega De Alatorre"
$ws.Range("B1540").Value = "Zozocolco De Hidalgo"
$ws.Range("B1557").Value = "Concepción Del Oro"
$ws.Range("B1565").Value = "Jiménez Del Teul"
$ws.Range("B1573").Value = "Nochistlán De Mejía"
$ws.Range("B1574").Value = "Noria De Ángeles"
$ws.Range("B1583").Value = "Teúl De González Ortega"
$ws.Range("B1584").Value = "Tlaltenango De Sánchez Román"
$ws.Range("B1588").Value = "Villa De Cos"

# --- Floating point bump (21 / 22483, recomputed one ULP higher) ---
$ws.Range("D4").Value = 0.0009340390517279724
$ws.Range("D449").Value = 0.0009340390517279724
$ws.Range("D627").Value = 0.0009340390517279724
$ws.Range("D1197").Value = 0.0009340390517279724
$ws.Range("D1252").Value = 0.0009340390517279724
$ws.Range("D1378").Value = 0.0009340390517279724
$ws.Range("D1379").Value = 0.0009340390517279724

# --- Drop the trailing footnote rows; dimension shrinks to A1:D1595 ---
$ws.Range("A1597:A1601").EntireRow.Delete()

